$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: remove H1 ("P368")
$ws.Range("H1").Value = $null

# Row 2: A2 "Gezinssamenstelling" -> "x"; remove H2 ("x")
$ws.Range("A2").Value = "x"
$ws.Range("H2").Value = $null

# Row 6: A6 "Eigen risico" -> "x"; remove H6 ("x")
$ws.Range("A6").Value = "x"
$ws.Range("H6").Value = $null

# Row 9: A9 "n/a" -> "x"; remove H9 ("x")
$ws.Range("A9").Value = "x"
$ws.Range("H9").Value = $null

# Row 14 was an empty gap row in the original sheet; populate B14 = "asd"
# without shifting any existing rows.
$ws.Range("B14").Value = "asd"
